# Handback report generation:
# - Status text "Ready for handoff" -> "Handed back: in sync with en-US" (Overview sheet,
#   columns E/F for both data rows; this is a shared string so editing one occurrence
#   updates every cell that referenced it).
# - zh-cn / de-de sheets: fill in "Latest Target File" (I), "Latest Handback File" (J)
#   and "Latest Handback DateTime" (K) for both data rows, and hyperlink the new
#   "Latest Target File" cells to the same source doc as column A.
# - Widen a few columns that now hold longer text/filenames.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: handoff status -> handback status for both rows / both locales
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Overview columns E (zh-cn) and F (de-de) grow to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Row 2 (099d4bef... file)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a744c2133b3d6bf4660e97168d79eb308851bfce/e2e/099d4bef-d7d9-445a-b17a-7624543dc0cf.md", "", "", "099d4bef-d7d9-445a-b17a-7624543dc0cf.md")
$zhcn.Range("J2").Value = "099d4bef-d7d9-445a-b17a-7624543dc0cf.bdaf9df29f3b9241e6815df1016e3f43641784e4.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-31 09:17:32"

# Row 3 (f15bc8f7... file)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a744c2133b3d6bf4660e97168d79eb308851bfce/e2e/f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f.md", "", "", "f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f.md")
$zhcn.Range("J3").Value = "f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f.325984f9a59ad99da4717e8dd1440fd382bd69d8.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-31 09:17:32"

# Column C (Status) grows to fit the longer status text; I/J grow to 40 chars for
# the new filenames, matching the Source File Name column's width.
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Row 2 (099d4bef... file)
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a744c2133b3d6bf4660e97168d79eb308851bfce/e2e/099d4bef-d7d9-445a-b17a-7624543dc0cf.md", "", "", "099d4bef-d7d9-445a-b17a-7624543dc0cf.md")
$dede.Range("J2").Value = "099d4bef-d7d9-445a-b17a-7624543dc0cf.bdaf9df29f3b9241e6815df1016e3f43641784e4.de-de.xlf"
$dede.Range("K2").Value = "2016-08-31 09:17:40"

# Row 3 (f15bc8f7... file)
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a744c2133b3d6bf4660e97168d79eb308851bfce/e2e/f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f.md", "", "", "f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f.md")
$dede.Range("J3").Value = "f15bc8f7-7e65-4ef3-9c56-d0dcd695fd4f.325984f9a59ad99da4717e8dd1440fd382bd69d8.de-de.xlf"
$dede.Range("K3").Value = "2016-08-31 09:17:40"

# Column C (Status) grows to fit the longer status text; I/J grow to 40 chars for
# the new filenames, matching the Source File Name column's width.
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
